$d = $word.ActiveDocument

# The revision-history table has a "description of change" cell whose
# paragraph reads "Dodane dodatne informacije u tekst te dodana slika".
# We append a further run containing " grafičkog prikaza" to that same
# paragraph (as its own run, not merged into the existing text), leaving
# the rest of the paragraph - including the two pre-existing runs -
# untouched.

$rng = $d.Content
$found = $rng.Find.Execute("Dodane dodatne informacije u tekst te dodana slika", `
                            $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate target paragraph text"
}

$para = $rng.Paragraphs(1)
$paraRange = $para.Range

$paraText = $paraRange.Text.TrimEnd([char]13, [char]7)
if ($paraText -ne "Dodane dodatne informacije u tekst te dodana slika") {
    throw "Unexpected paragraph text: [$paraText]"
}

# Rebuild the paragraph's OOXML exactly as it was, with one additional
# trailing run appended. The paragraph's own attributes (rsids, paraId,
# ...), its <w:pPr>, and the two existing <w:r> runs are reproduced
# byte-for-byte; only the new run is new.
$openTag = '<w:p w14:paraId="34CE0A41" w14:textId="5F2E4210" w:rsidR="00431A6E" w:rsidRPr="000833BA" w:rsidRDefault="00D50B1D" w:rsidP="00D50B1D">'
$pPr = '<w:pPr><w:tabs><w:tab w:val="left" w:pos="916"/><w:tab w:val="left" w:pos="1832"/><w:tab w:val="left" w:pos="2748"/><w:tab w:val="left" w:pos="3664"/><w:tab w:val="left" w:pos="4580"/><w:tab w:val="left" w:pos="5496"/><w:tab w:val="left" w:pos="6412"/><w:tab w:val="left" w:pos="7328"/><w:tab w:val="left" w:pos="8244"/><w:tab w:val="left" w:pos="9160"/><w:tab w:val="left" w:pos="10076"/><w:tab w:val="left" w:pos="10992"/><w:tab w:val="left" w:pos="11908"/><w:tab w:val="left" w:pos="12824"/><w:tab w:val="left" w:pos="13740"/><w:tab w:val="left" w:pos="14656"/></w:tabs><w:spacing w:after="60" w:line="240" w:lineRule="auto"/></w:pPr>'
$run1 = '<w:r w:rsidRPr="00D50B1D"><w:t>Dodane dodatne informacije u tekst</w:t></w:r>'
$run2 = '<w:r w:rsidR="007329CA"><w:t xml:space="preserve"> te dodana slika</w:t></w:r>'
$run3 = '<w:r><w:t xml:space="preserve"> grafičkog prikaza</w:t></w:r>'

$rebuilt = $openTag + $pPr + $run1 + $run2 + $run3 + "</w:p>"

$xml = @"
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
        <w:body>$rebuilt</w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

[void]$paraRange.InsertXML($xml)
